$d = $word.ActiveDocument

# The edit removes the explicit <w:contextualSpacing w:val="0"/> element
# from every paragraph's paragraph-properties (w:pPr) in the document.
# Word's object model does not expose a working ContextualSpacing
# property through this COM surface, so we rewrite each paragraph's
# OOXML (via Range.WordOpenXML / Range.InsertXML) with that element
# stripped out, leaving everything else about the paragraph untouched.

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $rng = $p.Range
    $xml = $rng.WordOpenXML
    if ($xml -like "*<w:contextualSpacing*") {
        $newXml = $xml -replace '<w:contextualSpacing[^/]*/>', ''
        [void]$rng.InsertXML($newXml)
    }
}
